$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "availability" column (G) marking each product's stock status.
$ws.Range("G1").Value = "availability"
$ws.Range("G2").Value = "Available"
$ws.Range("G3").Value = "Available"
$ws.Range("G4").Value = "Available"

# Match formatting of the existing header/data columns.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("F2").Copy()
$ws.Range("G2:G4").PasteSpecial(-4122)  # xlPasteFormats

$ws.Application.CutCopyMode = $false

# Keep the default print orientation explicit, matching a normal Excel resave.
$ws.PageSetup.Orientation = 1  # xlPortrait

# Leave the cursor on the newly added status cell, like the author likely did.
$ws.Range("G2").Select() | Out-Null
